$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark row 15 (LeetCode 503 - Next Greater Element I) as Done instead of On-going ---
$ws.Range("E2").Copy()
$ws.Range("E15").PasteSpecial(-4122)   # xlPasteFormats (brings style s="6")
$ws.Range("E15").Value = "Done"
$excel.CutCopyMode = 0

# --- Add three new Tree problems chosen by pedoe (rows 24-26) ---

# Row 24: Find Largest Value in Each Tree Row (515, Medium, On-going)
$ws.Range("A17:G17").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("F13").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A24").Value = 515
$ws.Range("B24").Value = "Find Largest Value in Each Tree Row"
$ws.Range("C24").Value = "Tree"
$ws.Range("D24").Value = "Pedoe"
$ws.Range("E24").Value = "On-going"
$ws.Range("F24").Value = "Medium"
$ws.Range("G24").Value = "Javascript"

# Row 25: Sum of Left Leaves (404, Easy, On-going)
$ws.Range("A17:G17").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A25").Value = 404
$ws.Range("B25").Value = "Sum of Left Leaves"
$ws.Range("C25").Value = "Tree"
$ws.Range("D25").Value = "Pedoe"
$ws.Range("E25").Value = "On-going"
$ws.Range("F25").Value = "Easy"
$ws.Range("G25").Value = "Javascript"

# Row 26: Symmetric Tree (101, Easy, On-going)
$ws.Range("A17:G17").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A26").Value = 101
$ws.Range("B26").Value = "Symmetric Tree"
$ws.Range("C26").Value = "Tree"
$ws.Range("D26").Value = "Pedoe"
$ws.Range("E26").Value = "On-going"
$ws.Range("F26").Value = "Easy"
$ws.Range("G26").Value = "Javascript"

# --- Update sheet view: selection on D27 (matches the author's last cursor position) ---
$ws.Range("D27").Select()
